# Weekly refresh of "Haba" (fava bean) price records for Feria Lagunitas de
# Puerto Montt: a couple of new daily samples were slotted in among the
# existing rows (pushing the later rows down one position each time) and
# three brand-new rows were appended at the bottom.
#
# Net effect on the sheet: rows 15-31 keep most of their static columns
# (Mercado ID/Mercado/Región/Codreg/Categoría/Variedad/Calidad/Unidad/
# Kg o Unidades/Clasificación) unchanged, but Fecha / Volumen / Precio
# mínimo / Precio máximo / Precio promedio ponderado / Origen / Precio $/Kg
# shift to the "next" row's former values (with two rows getting genuinely
# new data), and rows 32-34 are brand new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> Fecha(serial), Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg
$rows = @{
    15 = @(44428, 80,  18500, 18500, 18500, "Provincia de Limarí",    740)
    16 = @(44162, 140, 13000, 13000, 13000, "Región del Maule",       520)
    17 = @(44362, 60,  20000, 20000, 20000, "Provincia de Limarí",    800)
    18 = @(44421, 80,  17000, 17000, 17000, "Provincia de Limarí",    680)
    19 = @(44400, 80,  16500, 16500, 16500, "Provincia de Limarí",    660)
    20 = @(44390, 80,  16000, 16000, 16000, "Provincia de Limarí",    640)
    21 = @(44383, 80,  17000, 17000, 17000, "Provincia de Limarí",    680)
    22 = @(44169, 160, 13000, 14000, 13500, "Región de La Araucanía", 540)
    23 = @(44172, 40,  12000, 12000, 12000, "Región de La Araucanía", 480)
    24 = @(44214, 40,  25000, 25000, 25000, "Región de La Araucanía", 1000)
    25 = @(44435, 170, 18000, 19000, 18529, "Provincia de Limarí",    741)
    26 = @(44351, 30,  20000, 20000, 20000, "Región Metropolitana",   800)
    27 = @(44358, 60,  20000, 20000, 20000, "Región Metropolitana",   800)
    28 = @(44369, 70,  18000, 18000, 18000, "Provincia de Limarí",    720)
    29 = @(44161, 50,  11500, 11500, 11500, "Región del Maule",       460)
    30 = @(44159, 150, 11500, 11500, 11500, "Región del Maule",       460)
    31 = @(44160, 40,  11500, 11500, 11500, "Región del Maule",       460)
    32 = @(44376, 70,  17000, 17000, 17000, "Provincia de Limarí",    680)
    33 = @(44418, 90,  18000, 18000, 18000, "Provincia de Limarí",    720)
    34 = @(44432, 80,  18000, 18000, 18000, "Provincia de Limarí",    720)
}

foreach ($r in 15..34) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 1).Value  = 4
    $ws.Cells.Item($r, 2).Value  = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($r, 3).Value  = "Los Lagos"
    # Rows 32-34 are brand new, so the date column needs its number format
    # set explicitly (existing rows 15-31 already carry it from the template).
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 4).Value  = $vals[0]
    $ws.Cells.Item($r, 5).Value  = 10
    $ws.Cells.Item($r, 6).Value  = 100112026
    $ws.Cells.Item($r, 7).Value  = "Haba"
    $ws.Cells.Item($r, 8).Value  = "Sin especificar"
    $ws.Cells.Item($r, 9).Value  = "Primera"
    $ws.Cells.Item($r, 10).Value = $vals[1]
    $ws.Cells.Item($r, 11).Value = $vals[2]
    $ws.Cells.Item($r, 12).Value = $vals[3]
    $ws.Cells.Item($r, 13).Value = $vals[4]
    $ws.Cells.Item($r, 14).Value = "$/saco 25 kilos"
    $ws.Cells.Item($r, 15).Value = $vals[5]
    $ws.Cells.Item($r, 16).Value = $vals[6]
    $ws.Cells.Item($r, 17).Value = 25
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
